$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 (192 cell updates) ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 16:45:22'
$ws.Cells.Item(3, 1).Value = 'Total filas: 387'
$ws.Cells.Item(52, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(53, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(76, 3).Value = '10_OLMOS'
$ws.Cells.Item(77, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(124, 1).Value = '08:39:56'
$ws.Cells.Item(124, 3).Value = '215C_EL PATO'
$ws.Cells.Item(124, 4).Value = 62
$ws.Cells.Item(125, 1).Value = '09:38:04'
$ws.Cells.Item(125, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(125, 4).Value = 3
$ws.Cells.Item(126, 3).Value = '14_ABASTO'
$ws.Cells.Item(174, 1).Value = '09:38:04'
$ws.Cells.Item(174, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(174, 4).Value = 107
$ws.Cells.Item(175, 1).Value = '11:23:54'
$ws.Cells.Item(175, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(175, 4).Value = 2
$ws.Cells.Item(199, 1).Value = '10:28:12'
$ws.Cells.Item(199, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(199, 4).Value = 98
$ws.Cells.Item(200, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(201, 1).Value = '10:57:58'
$ws.Cells.Item(201, 3).Value = '14_ABASTO'
$ws.Cells.Item(201, 4).Value = 69
$ws.Cells.Item(212, 3).Value = '215A_EL PATO'
$ws.Cells.Item(213, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(224, 1).Value = '11:51:05'
$ws.Cells.Item(224, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(224, 4).Value = 46
$ws.Cells.Item(225, 1).Value = '10:57:58'
$ws.Cells.Item(225, 3).Value = '17_179 Y 38'
$ws.Cells.Item(225, 4).Value = 100
$ws.Cells.Item(226, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(307, 1).Value = '13:55:06'
$ws.Cells.Item(307, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(307, 4).Value = 84
$ws.Cells.Item(308, 1).Value = '15:19:52'
$ws.Cells.Item(308, 3).Value = '10_OLMOS'
$ws.Cells.Item(308, 4).Value = 0
$ws.Cells.Item(330, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(332, 3).Value = '10_OLMOS'
$ws.Cells.Item(357, 1).Value = '16:45:22'
$ws.Cells.Item(357, 2).Value = '16:45'
$ws.Cells.Item(357, 3).Value = '14_ABASTO'
$ws.Cells.Item(357, 4).Value = 0
$ws.Cells.Item(358, 1).Value = '14:53:55'
$ws.Cells.Item(358, 2).Value = '16:48'
$ws.Cells.Item(358, 3).Value = '15_ABASTO'
$ws.Cells.Item(358, 4).Value = 115
$ws.Cells.Item(359, 1).Value = '15:51:40'
$ws.Cells.Item(359, 2).Value = '16:50'
$ws.Cells.Item(359, 3).Value = '14_ABASTO'
$ws.Cells.Item(359, 4).Value = 59
$ws.Cells.Item(360, 1).Value = '15:19:52'
$ws.Cells.Item(360, 3).Value = '17_179 Y 38'
$ws.Cells.Item(360, 4).Value = 97
$ws.Cells.Item(361, 1).Value = '16:14:52'
$ws.Cells.Item(361, 2).Value = '16:56'
$ws.Cells.Item(361, 4).Value = 42
$ws.Cells.Item(362, 1).Value = '16:32:38'
$ws.Cells.Item(362, 2).Value = '16:57'
$ws.Cells.Item(362, 3).Value = '10_OLMOS'
$ws.Cells.Item(362, 4).Value = 25
$ws.Cells.Item(363, 1).Value = '16:14:52'
$ws.Cells.Item(363, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(363, 4).Value = 50
$ws.Cells.Item(364, 1).Value = '15:51:40'
$ws.Cells.Item(364, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(364, 4).Value = 73
$ws.Cells.Item(365, 1).Value = '15:19:52'
$ws.Cells.Item(365, 2).Value = '17:04'
$ws.Cells.Item(365, 3).Value = '215A_EL PATO'
$ws.Cells.Item(365, 4).Value = 105
$ws.Cells.Item(366, 1).Value = '16:45:22'
$ws.Cells.Item(366, 2).Value = '17:05'
$ws.Cells.Item(366, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(366, 4).Value = 20
$ws.Cells.Item(367, 2).Value = '17:10'
$ws.Cells.Item(367, 3).Value = '10_OLMOS'
$ws.Cells.Item(367, 4).Value = 38
$ws.Cells.Item(368, 1).Value = '16:45:22'
$ws.Cells.Item(368, 2).Value = '17:16'
$ws.Cells.Item(368, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(368, 4).Value = 31
$ws.Cells.Item(369, 1).Value = '16:32:38'
$ws.Cells.Item(369, 2).Value = '17:20'
$ws.Cells.Item(369, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(369, 4).Value = 48
$ws.Cells.Item(370, 1).Value = '16:14:52'
$ws.Cells.Item(370, 2).Value = '17:20'
$ws.Cells.Item(370, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(370, 4).Value = 66
$ws.Cells.Item(371, 1).Value = '16:45:22'
$ws.Cells.Item(371, 2).Value = '17:21'
$ws.Cells.Item(371, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(371, 4).Value = 36
$ws.Cells.Item(372, 1).Value = '15:51:40'
$ws.Cells.Item(372, 2).Value = '17:21'
$ws.Cells.Item(372, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(372, 4).Value = 90
$ws.Cells.Item(373, 1).Value = '15:51:40'
$ws.Cells.Item(373, 2).Value = '17:24'
$ws.Cells.Item(373, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(373, 4).Value = 93
$ws.Cells.Item(374, 2).Value = '17:28'
$ws.Cells.Item(374, 3).Value = '14_ABASTO'
$ws.Cells.Item(374, 4).Value = 97
$ws.Cells.Item(375, 1).Value = '16:45:22'
$ws.Cells.Item(375, 2).Value = '17:31'
$ws.Cells.Item(375, 3).Value = '15_ABASTO'
$ws.Cells.Item(375, 4).Value = 46
$ws.Cells.Item(376, 1).Value = '16:32:38'
$ws.Cells.Item(376, 2).Value = '17:32'
$ws.Cells.Item(376, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(376, 4).Value = 60
$ws.Cells.Item(377, 1).Value = '16:14:52'
$ws.Cells.Item(377, 2).Value = '17:35'
$ws.Cells.Item(377, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(377, 4).Value = 81
$ws.Cells.Item(378, 1).Value = '16:32:38'
$ws.Cells.Item(378, 2).Value = '17:35'
$ws.Cells.Item(378, 3).Value = '15_ABASTO'
$ws.Cells.Item(378, 4).Value = 63
$ws.Cells.Item(379, 1).Value = '16:45:22'
$ws.Cells.Item(379, 2).Value = '17:35'
$ws.Cells.Item(379, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(379, 4).Value = 50
$ws.Cells.Item(380, 1).Value = '15:51:40'
$ws.Cells.Item(380, 2).Value = '17:36'
$ws.Cells.Item(380, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(380, 4).Value = 105
$ws.Cells.Item(381, 1).Value = '16:45:22'
$ws.Cells.Item(381, 2).Value = '17:37'
$ws.Cells.Item(381, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(381, 4).Value = 52
$ws.Cells.Item(382, 1).Value = '15:51:40'
$ws.Cells.Item(382, 2).Value = '17:38'
$ws.Cells.Item(382, 3).Value = '17_ROMERO'
$ws.Cells.Item(382, 4).Value = 107
$ws.Cells.Item(382, 5).Value = 'LP1912'
$ws.Cells.Item(383, 1).Value = '15:51:40'
$ws.Cells.Item(383, 2).Value = '17:40'
$ws.Cells.Item(383, 3).Value = '215B_EL PATO'
$ws.Cells.Item(383, 4).Value = 109
$ws.Cells.Item(383, 5).Value = 'LP1912'
$ws.Cells.Item(384, 1).Value = '16:45:22'
$ws.Cells.Item(384, 2).Value = '17:41'
$ws.Cells.Item(384, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(384, 4).Value = 56
$ws.Cells.Item(384, 5).Value = 'LP1912'
$ws.Cells.Item(385, 1).Value = '16:45:22'
$ws.Cells.Item(385, 2).Value = '17:45'
$ws.Cells.Item(385, 3).Value = '15_ABASTO'
$ws.Cells.Item(385, 4).Value = 60
$ws.Cells.Item(385, 5).Value = 'LP1912'
$ws.Cells.Item(386, 1).Value = '15:51:40'
$ws.Cells.Item(386, 2).Value = '17:50'
$ws.Cells.Item(386, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(386, 4).Value = 119
$ws.Cells.Item(386, 5).Value = 'LP1912'
$ws.Cells.Item(387, 1).Value = '16:14:52'
$ws.Cells.Item(387, 2).Value = '17:52'
$ws.Cells.Item(387, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(387, 4).Value = 98
$ws.Cells.Item(387, 5).Value = 'LP1912'
$ws.Cells.Item(388, 1).Value = '16:14:52'
$ws.Cells.Item(388, 2).Value = '18:04'
$ws.Cells.Item(388, 3).Value = '17_ROMERO'
$ws.Cells.Item(388, 4).Value = 110
$ws.Cells.Item(388, 5).Value = 'LP1912'
$ws.Cells.Item(389, 1).Value = '16:32:38'
$ws.Cells.Item(389, 2).Value = '18:21'
$ws.Cells.Item(389, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(389, 4).Value = 109
$ws.Cells.Item(389, 5).Value = 'LP1912'
$ws.Cells.Item(390, 1).Value = '16:32:38'
$ws.Cells.Item(390, 2).Value = '18:27'
$ws.Cells.Item(390, 3).Value = '215C_EL PATO'
$ws.Cells.Item(390, 4).Value = 115
$ws.Cells.Item(390, 5).Value = 'LP1912'
$ws.Cells.Item(391, 1).Value = '16:45:22'
$ws.Cells.Item(391, 2).Value = '18:28'
$ws.Cells.Item(391, 3).Value = '215C_EL PATO'
$ws.Cells.Item(391, 4).Value = 103
$ws.Cells.Item(391, 5).Value = 'LP1912'
$ws.Cells.Item(392, 1).Value = '16:45:22'
$ws.Cells.Item(392, 2).Value = '18:32'
$ws.Cells.Item(392, 3).Value = '11X44_ETCHEVERRY'
$ws.Cells.Item(392, 4).Value = 107
$ws.Cells.Item(392, 5).Value = 'LP1912'

# ---- Sheet: LP1912-215 (7 cell updates) ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 16:45:22'
$ws.Cells.Item(3, 1).Value = 'Total filas: 40'
$ws.Cells.Item(45, 1).Value = '16:45:22'
$ws.Cells.Item(45, 2).Value = '18:28'
$ws.Cells.Item(45, 3).Value = '215C_EL PATO'
$ws.Cells.Item(45, 4).Value = 103
$ws.Cells.Item(45, 5).Value = 'LP1912'

# ---- Sheet: 6203-6173 (7 cell updates) ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 16:45:22'
$ws.Cells.Item(3, 1).Value = 'Total filas: 52'
$ws.Cells.Item(57, 1).Value = '16:45:22'
$ws.Cells.Item(57, 2).Value = '18:04'
$ws.Cells.Item(57, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(57, 4).Value = 79
$ws.Cells.Item(57, 5).Value = 'L6203'

Write-Output "done"